$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A8").Value = "TV4"
$ws.Range("B8").Value = "S3"
$ws.Range("C8").Value = "PYR"
$ws.Range("D8").Value = "MC1"
$ws.Range("E8").Value = 5
$ws.Range("F8").Value = "MC3"
$ws.Range("G8").Value = 1
$ws.Range("H8").Value = 45401
$ws.Range("H7").Copy()
$ws.Range("H8").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("I8").Value = "Port went bad"

$ws.Range("H8").Select()
